$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J: "Maple_Died" indicator
$ws.Range("J1").Value = "Maple_Died"

# Rows where the Maple died ("Y")
$yRows = @(25, 60, 62, 185, 325, 525)
foreach ($r in $yRows) {
    $ws.Cells.Item($r, 10).Value = "Y"
}

# New family data: fill the previously-empty Sampler column (I) for rows
# 310-558 with "Nikita"
for ($r = 310; $r -le 558; $r++) {
    $ws.Cells.Item($r, 9).Value = "Nikita"
}

# Restore the workbook's final on-screen selection
$ws.Range("K11").Select()
